$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that currently carry a password-ish value in column F; that value
# needs to move one column right (to G) and column F gets a new "*" mask
# cell (unstyled) in its place.
$rows = @(7, 8, 9, 10, 13, 14, 15)

foreach ($r in $rows) {
    $src = $ws.Cells.Item($r, 6)   # F{r}
    $dst = $ws.Cells.Item($r, 7)   # G{r}
    $src.Copy($dst)
}

# Re-home the three hyperlinks that lived on F9/F10/F14 onto the cells that
# now hold the real values (G9/G10/G14). The engine's Hyperlink.Delete is a
# no-op, but Range.Hyperlinks.Delete() clears every hyperlink on the sheet,
# so clear once and re-add in the original order to keep rId1/rId2/rId3
# aligned with G9/G10/G14.
$ws.Range("A1").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Cells.Item(9, 7), "about:blank")
$ws.Hyperlinks.Add($ws.Cells.Item(10, 7), "about:blank")
$ws.Hyperlinks.Add($ws.Cells.Item(14, 7), "about:blank")

foreach ($r in $rows) {
    $f = $ws.Cells.Item($r, 6)
    $f.ClearFormats()
    $f.Value = "*"
}

# Column widths: E keeps its own width (now best-fit), F takes the old
# E:F width, G (the shifted password column) gets a new, wider column.
$ws.Columns.Item(5).ColumnWidth = 24.14
$ws.Columns.Item(6).ColumnWidth = 15
$ws.Columns.Item(7).ColumnWidth = 18.29

# Sheet view: scrolled down with A2:A34 selected (A34 active).
$ws.Application.GoTo($ws.Range("A10"))
$ws.Range("A2:A34").Select()
$ws.Cells.Item(34, 1).Activate()
